$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("116:116").Insert()

$ws.Cells.Item(116, 1).Value = 11
$ws.Cells.Item(116, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(116, 3).Value = "Bíobío"
$ws.Cells.Item(116, 4).Value = 44484
$ws.Cells.Item(116, 5).Value = 8
$ws.Cells.Item(116, 6).Value = 100114014
$ws.Cells.Item(116, 7).Value = "Betarraga"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 1500
$ws.Cells.Item(116, 11).Value = 600
$ws.Cells.Item(116, 12).Value = 700
$ws.Cells.Item(116, 13).Value = 647
$ws.Cells.Item(116, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(116, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 16).Value = 129
$ws.Cells.Item(116, 17).Value = 5
$ws.Cells.Item(116, 18).Value = "Hortaliza"
